$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'24.807.60"
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = "'1.708.05"
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("D4").Value = "'0.9982"
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").Value = "'317.91"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").Value = "'0.9973"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = "'0.3923"
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").Value = "'0.4069"
$ws.Range("E8").Value = '  +0.68%  '
$ws.Range("D9").Value = "'1.499"
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("B10").Value = 'BinanceUSD'
$ws.Range("C10").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D10").Value = "'0.9966"
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").Value = "'53.41"
$ws.Range("E11").Value = '  +1.64%  '
$ws.Range("D12").Value = "'0.08828"
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").Value = "'26.52"
$ws.Range("E13").Value = '  +12.40%  '
$ws.Range("D14").Value = "'7.514"
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = "'8.134"
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("D16").Value = "'0.00001362"
$ws.Range("E16").Value = '  +3.11%  '
$ws.Range("D17").Value = "'1.704.19"
$ws.Range("E17").Value = '  +0.84%  '
$ws.Range("D18").Value = "'97.69"
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("D19").Value = "'0.07195"
$ws.Range("E19").Value = '  +1.79%  '
$ws.Range("D20").Value = "'20.61"
$ws.Range("E20").Value = '  +4.16%  '
$ws.Range("D21").Value = "'7.312"
$ws.Range("E21").Value = '  +3.09%  '
$ws.Range("D22").Value = "'0.9972"
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("D23").Value = "'14.43"
$ws.Range("E23").Value = '  -2.11%  '
$ws.Range("D24").Value = "'24.821.35"
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").Value = "'3.027"
$ws.Range("E25").Value = '  -3.42%  '
$ws.Range("D26").Value = "'2.335"
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("D27").Value = "'23.05"
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").Value = "'167.75"
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("D29").Value = "'5.959"
$ws.Range("E29").Value = '  +15.83%  '
$ws.Range("D30").Value = "'8.578"
$ws.Range("E30").Value = '  -6.78%  '
$ws.Range("D31").Value = "'145.07"
$ws.Range("E31").Value = '  +6.87%  '
$ws.Range("D32").Value = "'1.893.50"
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("B33").Value = 'WEMIXTOKEN'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = "'2.192"
$ws.Range("E33").Value = '  +11.49%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = "'0.08843"
$ws.Range("E34").Value = '  -2.06%  '
$ws.Range("D35").Value = "'1.060"
$ws.Range("E35").Value = '  -0.86%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = "'7.247"
$ws.Range("E36").Value = '  -9.84%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = "'0.03142"
$ws.Range("E37").Value = '  +5.91%  '
$ws.Range("D38").Value = "'0.2823"
$ws.Range("E38").Value = '  +1.72%  '
$ws.Range("D39").Value = "'10.97"
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("D40").Value = "'0.8474"
$ws.Range("E40").Value = '  +9.21%  '
$ws.Range("D41").Value = "'0.09247"
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = "'14.17"
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("D43").Value = "'1.476"
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("D44").Value = "'17.84"
$ws.Range("E44").Value = '  +11.20%  '
$ws.Range("D45").Value = "'2.710"
$ws.Range("E45").Value = '  +4.38%  '
$ws.Range("D46").Value = "'0.7499"
$ws.Range("E46").Value = '  +4.00%  '
$ws.Range("D47").Value = "'4.286"
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("D48").Value = "'1.395"
$ws.Range("E48").Value = '  +3.11%  '
$ws.Range("D49").Value = "'0.9967"
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").Value = "'140.75"
$ws.Range("E50").Value = '  +0.90%  '
$ws.Range("D51").Value = "'0.08280"
$ws.Range("E51").Value = '  +3.65%  '
